$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking text (e.g. "40.119.70", "1.00", "0.100")
# Excel auto-converts these to numbers on assignment, which both changes the
# stored type and loses formatting (trailing zeros, multi-dot groupings).
# Force text storage via NumberFormat="@", then restore the original style so
# the cell style index is unchanged in the saved file.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "40.119.70"
$ws.Range("E2").Value = "  +1.84%  "

Set-TextValue $ws.Range("D3") "2.236.96"
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "293.43"
$ws.Range("E5").Value = "  -1.16%  "

Set-TextValue $ws.Range("D6") "87.25"
$ws.Range("E6").Value = "  +5.02%  "

$ws.Range("E7").Value = "  +1.18%  "

$ws.Range("E8").Value = "  -0.03%  "

Set-TextValue $ws.Range("D9") "0.476"
$ws.Range("E9").Value = "  +1.51%  "

Set-TextValue $ws.Range("D10") "31.30"
$ws.Range("E10").Value = "  +7.37%  "

Set-TextValue $ws.Range("D11") "0.0791"
$ws.Range("E11").Value = "  +2.11%  "

Set-TextValue $ws.Range("D12") "46.99"
$ws.Range("E12").Value = "  -1.76%  "

$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("E14").Value = "  +1.76%  "

Set-TextValue $ws.Range("D15") "2.583.97"
$ws.Range("E15").Value = "  +0.89%  "

Set-TextValue $ws.Range("D16") "14.11"
$ws.Range("E16").Value = "  +0.02%  "

Set-TextValue $ws.Range("D17") "2.239.49"
$ws.Range("E17").Value = "  +1.49%  "

$ws.Range("E18").Value = "  +2.61%  "

Set-TextValue $ws.Range("D19") "40.028.53"
$ws.Range("E19").Value = "  +1.76%  "

Set-TextValue $ws.Range("D20") "0.0₃0891"
$ws.Range("E20").Value = "  +1.86%  "

Set-TextValue $ws.Range("D21") "11.26"
$ws.Range("E21").Value = "  +9.47%  "

Set-TextValue $ws.Range("D22") "5.84"
$ws.Range("E22").Value = "  +2.14%  "

Set-TextValue $ws.Range("D23") "65.83"
$ws.Range("E23").Value = "  +1.20%  "

Set-TextValue $ws.Range("D24") "236.32"
$ws.Range("E24").Value = "  +3.76%  "

Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  +3.11%  "

Set-TextValue $ws.Range("D27") "1.86"
$ws.Range("E27").Value = "  +2.27%  "

Set-TextValue $ws.Range("D28") "23.01"
$ws.Range("E28").Value = "  +1.87%  "

$ws.Range("E29").Value = "  +2.62%  "

Set-TextValue $ws.Range("D30") "9.35"
$ws.Range("E30").Value = "  +2.41%  "

Set-TextValue $ws.Range("D31") "33.44"
$ws.Range("E31").Value = "  +4.20%  "

Set-TextValue $ws.Range("D32") "151.57"
$ws.Range("E32").Value = "  +1.16%  "

Set-TextValue $ws.Range("D33") "0.998"
$ws.Range("E33").Value = "  -0.25%  "

Set-TextValue $ws.Range("D34") "4.94"
$ws.Range("E34").Value = "  +1.82%  "

Set-TextValue $ws.Range("D35") "0.0722"
$ws.Range("E35").Value = "  +3.97%  "

$ws.Range("E36").Value = "  +2.04%  "

$ws.Range("E37").Value = "  +7.14%  "

Set-TextValue $ws.Range("D38") "2.83"
$ws.Range("E38").Value = "  +6.96%  "

$ws.Range("E39").Value = "  +2.18%  "

Set-TextValue $ws.Range("D40") "0.100"
$ws.Range("E40").Value = "  +3.68%  "

$ws.Range("E41").Value = "  +4.81%  "

$ws.Range("E42").Value = "  +5.06%  "

Set-TextValue $ws.Range("D43") "2.063.44"
$ws.Range("E43").Value = "  +8.10%  "

Set-TextValue $ws.Range("D44") "18.19"
$ws.Range("E44").Value = "  +13.25%  "

$ws.Range("E45").Value = "  +4.21%  "

$ws.Range("E46").Value = "  +4.04%  "

$ws.Range("E47").Value = "  +9.57%  "

$ws.Range("E48").Value = "  -0.38%  "

$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D49") "72.30"
$ws.Range("E49").Value = "  +2.37%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D50") "2.443.97"
$ws.Range("E50").Value = "  +0.51%  "

Set-TextValue $ws.Range("D51") "89.47"
$ws.Range("E51").Value = "  +2.81%  "

